# Update column F (dSF) values for specific rows based on repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 8
    3  = -1
    4  = -2
    5  = -8
    8  = 8
    13 = -4
    14 = -2
    15 = -3
    18 = -3
    19 = 0
    22 = 5
    23 = -2
    28 = -3
    29 = 3
    32 = -1
    33 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
